$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.217.05'
$ws.Range('E2').Value = '  +0.59%  '
$ws.Range('D3').Value = '2.602.59'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '540.53'
$ws.Range('E5').Value = '  +3.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.49'
$ws.Range('E6').Value = '  +1.23%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.565'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.43'
$ws.Range('E9').Value = '  -1.67%  '
$ws.Range('E10').Value = '  +1.95%  '
$ws.Range('E11').Value = '  +1.47%  '
$ws.Range('E12').Value = '  +0.72%  '
$ws.Range('D13').Value = '3.060.01'
$ws.Range('D14').Value = '59.162.87'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.57'
$ws.Range('E15').Value = '  +0.44%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.659.04'
$ws.Range('E16').Value = '  +2.13%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000134'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '341.03'
$ws.Range('E18').Value = '  +0.63%  '
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.12'
$ws.Range('E21').Value = '  -1.94%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.62'
$ws.Range('E23').Value = '  +2.29%  '
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('E25').Value = '  -1.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.19'
$ws.Range('E27').Value = '  +1.89%  '
$ws.Range('D28').Value = '0.0₃0751'
$ws.Range('E28').Value = '  +3.77%  '
$ws.Range('E29').Value = '  +0.06%  '
$ws.Range('E30').Value = '  +7.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.81'
$ws.Range('E31').Value = '  -2.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.72'
$ws.Range('E32').Value = '  -0.30%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '149.67'
$ws.Range('E33').Value = '  +0.45%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.98'
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.11'
$ws.Range('E35').Value = '  -0.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '37.14'
$ws.Range('E36').Value = '  +2.36%  '
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.833'
$ws.Range('E38').Value = '  +0.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.826'
$ws.Range('E39').Value = '  +0.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.56'
$ws.Range('E40').Value = '  +1.62%  '
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '274.86'
$ws.Range('E42').Value = '  +0.30%  '
$ws.Range('E43').Value = '  +1.30%  '
$ws.Range('E44').Value = '  -0.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0955'
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('E46').Value = '  +0.33%  '
$ws.Range('D47').Value = '1.954.95'
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.52'
$ws.Range('E48').Value = '  +3.35%  '
$ws.Range('E49').Value = '  +1.26%  '
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '111.31'
$ws.Range('E51').Value = '  -0.22%  '
